$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1118.9615
$ws.Range("J17").Value = 923.9804
$ws.Range("L17").Value = 2771.9412
$ws.Range("N17").Value = -3107.9412
$ws.Range("H18").Value = 14728
$ws.Range("I18").Value = 2620
$ws.Range("J18").Value = 20109.334
$ws.Range("K18").Value = 2620
$ws.Range("L18").Value = 20109.334
$ws.Range("M18").Value = -2336
$ws.Range("N18").Value = -20677.334
$ws.Range("H127").Value = 5058.25
$ws.Range("I127").Value = 6244.3335
$ws.Range("K127").Value = 18733.0005
$ws.Range("M127").Value = -13773.0005
$ws.Range("H131").Value = 2265
$ws.Range("I131").Value = 1047.1428
$ws.Range("J131").Value = 3685.8333
$ws.Range("K131").Value = 3141.4284
$ws.Range("L131").Value = 11057.4999
$ws.Range("M131").Value = 1898.5716
$ws.Range("N131").Value = -21137.4999
$ws.Range("H137").Value = 1623.5333
$ws.Range("I137").Value = 1244.5454
$ws.Range("K137").Value = 3733.6362
$ws.Range("M137").Value = -1183.6362
$ws.Range("H138").Value = 1773.8334
$ws.Range("J138").Value = 2211.5557
$ws.Range("L138").Value = 6634.6671
$ws.Range("N138").Value = -16914.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 6879.6
$ws.Range("I80").Value = 46.5
$ws.Range("J80").Value = 7930.846
$ws.Range("K80").Value = 46.5
$ws.Range("L80").Value = 7930.846
$ws.Range("M80").Value = 951.5
$ws.Range("N80").Value = -9926.846
$ws.Range("H82").Value = 28750
$ws.Range("H83").Value = 6879.6
$ws.Range("I83").Value = 46.5
$ws.Range("J83").Value = 7930.846
$ws.Range("K83").Value = 232.5
$ws.Range("L83").Value = 39654.23
$ws.Range("M83").Value = 4759.5
$ws.Range("N83").Value = -49638.23
$ws.Range("H85").Value = 28750

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2026.5
$ws.Range("I134").Value = 1239.2
$ws.Range("K134").Value = 3717.6
$ws.Range("M134").Value = -1182.6
$ws.Range("H135").Value = 27857.143
$ws.Range("I135").Value = 15000
$ws.Range("J135").Value = 30000
$ws.Range("K135").Value = 15000
$ws.Range("L135").Value = 30000
$ws.Range("M135").Value = -9930
$ws.Range("N135").Value = -40140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1999
$ws.Range("I3").Value = 1999
$ws.Range("K3").Value = 5997
$ws.Range("M3").Value = -5885
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H56").Value = 8068.4546
$ws.Range("I56").Value = 8068.4546
$ws.Range("K56").Value = 8068.4546
$ws.Range("M56").Value = -7538.4546
$ws.Range("H107").Value = 703.55
$ws.Range("J107").Value = 703.55
$ws.Range("L107").Value = 2110.65
$ws.Range("N107").Value = -5950.65
$ws.Range("H114").Value = 47619520
$ws.Range("J114").Value = 142857140
$ws.Range("L114").Value = 428571420
$ws.Range("N114").Value = -428577928
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H131").Value = 11583.815
$ws.Range("I131").Value = 532
$ws.Range("J131").Value = 12308.524
$ws.Range("K131").Value = 1596
$ws.Range("L131").Value = 36925.572
$ws.Range("M131").Value = 3444
$ws.Range("N131").Value = -47005.572
$ws.Range("H133").Value = 3629.9092
$ws.Range("I133").Value = 1232.25
$ws.Range("K133").Value = 3696.75
$ws.Range("M133").Value = 1363.25
$ws.Range("H134").Value = 1440.2609
$ws.Range("I134").Value = 1196.6666
$ws.Range("K134").Value = 3589.9998
$ws.Range("M134").Value = 1480.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5691.923
$ws.Range("I70").Value = 6059.6
$ws.Range("K70").Value = 6059.6
$ws.Range("M70").Value = -5789.6
$ws.Range("H73").Value = 5691.923
$ws.Range("I73").Value = 6059.6
$ws.Range("K73").Value = 6059.6
$ws.Range("M73").Value = -5123.6
$ws.Range("H80").Value = 2265.4443
$ws.Range("I80").Value = 1863
$ws.Range("J80").Value = 2466.6667
$ws.Range("K80").Value = 1863
$ws.Range("L80").Value = 2466.6667
$ws.Range("M80").Value = -865
$ws.Range("N80").Value = -4462.6667
$ws.Range("H83").Value = 2265.4443
$ws.Range("I83").Value = 1863
$ws.Range("J83").Value = 2466.6667
$ws.Range("K83").Value = 9315
$ws.Range("L83").Value = 12333.3335
$ws.Range("M83").Value = -4323
$ws.Range("N83").Value = -22317.3335
$ws.Range("H107").Value = 516.6667
$ws.Range("I107").Value = 420
$ws.Range("K107").Value = 420
$ws.Range("M107").Value = 1500
$ws.Range("H110").Value = 99995
$ws.Range("J110").Value = 99995
$ws.Range("L110").Value = 99995
$ws.Range("N110").Value = -108175

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -705
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -893
$ws.Range("N27").ClearContents()
$ws.Range("H46").Value = 1700
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1700
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 1700
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -2076
$ws.Range("H68").Value = 5282.5713
$ws.Range("I68").Value = 5494.75
$ws.Range("K68").Value = 5494.75
$ws.Range("M68").Value = -4745.75
$ws.Range("H71").Value = 5282.5713
$ws.Range("I71").Value = 5494.75
$ws.Range("K71").Value = 27473.75
$ws.Range("M71").Value = -23729.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 70007
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H113").Value = 616.3182
$ws.Range("I113").Value = 411.70587
$ws.Range("J113").Value = 1312
$ws.Range("K113").Value = 1235.11761
$ws.Range("L113").Value = 3936
$ws.Range("M113").Value = 934.88239
$ws.Range("N113").Value = -8276
